$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix SOM (organic matter) and summary calculations ---
# B4: Sample_ID 19 -> value corrected from 100 to 102
$ws.Range("B4").Value = 102

# B12/B13: Organic_matter / Bulk_Density rows for sample 28/29 were placeholders (0);
# corrected to their real measured values
$ws.Range("B12").Value = 1.2
$ws.Range("B13").Value = 1.3

# B14: Layer_depth for sample 30 corrected from 50 to 30
$ws.Range("B14").Value = 30

# --- Restore the view scroll position / selection as last saved by the author ---
$ws.Activate()
$ws.Range("B14").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
